$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header volume/issue number and date range text
$a8 = $ws.Range("A8").Text
$idxNum = $a8.IndexOf("24") + 1
$ws.Range("A8").Characters($idxNum, 2).Text = "25"

$c9 = $ws.Range("C9").Text
$idxD1 = $c9.IndexOf("6/12/2023") + 1
$ws.Range("C9").Characters($idxD1, 9).Text = "6/19/2023"
$c9b = $ws.Range("C9").Text
$idxD2 = $c9b.IndexOf("6/18/2023") + 1
$ws.Range("C9").Characters($idxD2, 9).Text = "6/25/2023"

# Update crime statistics table cells
$ws.Range("D15").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = 1
$ws.Range("E15").Value = 0
$ws.Range("D15").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = -50
$ws.Range("I15").Value = 5
$ws.Range("J15").Value = 6
$ws.Range("K15").Value = -16.666666666666
$ws.Range("M15").Value = -37.5
$ws.Range("N15").Value = -64.285714285714
$ws.Range("C16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 140
$ws.Range("I16").Value = 48
$ws.Range("J16").Value = 53
$ws.Range("K16").Value = -9.433962264150
$ws.Range("L16").Value = 54.838709677419
$ws.Range("M16").Value = -54.716981132075
$ws.Range("N16").Value = -85.276073619631
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 166.666666666667
$ws.Range("F17").Value = 22
$ws.Range("H17").Value = 69.230769230769
$ws.Range("I17").Value = 83
$ws.Range("J17").Value = 80
$ws.Range("K17").Value = 3.75
$ws.Range("L17").Value = 13.698630136986
$ws.Range("M17").Value = 38.333333333333
$ws.Range("N17").Value = -44.295302013422
$ws.Range("D15").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C18").Value = 4
$ws.Range("A14").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Formula = "=""0"""
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("A14").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E18").Formula = "=""***.*"""
$ws.Range("E18").Copy()
$ws.Range("E18").PasteSpecial(-4163)
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 100
$ws.Range("I18").Value = 50
$ws.Range("K18").Value = 6.382978723404
$ws.Range("L18").Value = -1.960784313725
$ws.Range("M18").Value = -57.627118644067
$ws.Range("N18").Value = -91.694352159468
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 33.333333333333
$ws.Range("F19").Value = 39
$ws.Range("G19").Value = 48
$ws.Range("H19").Value = -18.75
$ws.Range("I19").Value = 286
$ws.Range("J19").Value = 233
$ws.Range("K19").Value = 22.746781115879
$ws.Range("L19").Value = 90.666666666666
$ws.Range("M19").Value = 25.438596491228
$ws.Range("N19").Value = -13.855421686747
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = 25
$ws.Range("I20").Value = 52
$ws.Range("J20").Value = 54
$ws.Range("K20").Value = -3.703703703703
$ws.Range("L20").Value = 23.809523809523
$ws.Range("M20").Value = -34.177215189873
$ws.Range("N20").Value = -96.226415094339
$ws.Range("C21").Value = 30
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = 76.470588235294
$ws.Range("F21").Value = 92
$ws.Range("G21").Value = 80
$ws.Range("H21").Value = 15
$ws.Range("I21").Value = 526
$ws.Range("J21").Value = 473
$ws.Range("K21").Value = 11.205073995771
$ws.Range("L21").Value = 49.008498583569
$ws.Range("M21").Value = -12.186978297161
$ws.Range("N21").Value = -81.287797936677
$ws.Range("D15").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C23").Value = 1
$ws.Range("A14").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("D23").Formula = "=""0"""
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("A14").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("E23").Formula = "=""***.*"""
$ws.Range("E23").Copy()
$ws.Range("E23").PasteSpecial(-4163)
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = -60
$ws.Range("I23").Value = 8
$ws.Range("K23").Value = -57.894736842105
$ws.Range("L23").Value = 60
$ws.Range("M23").Value = -50
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 32
$ws.Range("E24").Value = -34.375
$ws.Range("F24").Value = 91
$ws.Range("G24").Value = 92
$ws.Range("H24").Value = -1.086956521739
$ws.Range("I24").Value = 576
$ws.Range("J24").Value = 520
$ws.Range("K24").Value = 10.769230769230
$ws.Range("L24").Value = 57.377049180327
$ws.Range("M24").Value = 32.413793103448
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 60
$ws.Range("F25").Value = 28
$ws.Range("G25").Value = 22
$ws.Range("H25").Value = 27.272727272727
$ws.Range("I25").Value = 128
$ws.Range("J25").Value = 107
$ws.Range("K25").Value = 19.626168224299
$ws.Range("L25").Value = 10.344827586206
$ws.Range("M25").Value = -23.353293413173
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = -33.333333333333
$ws.Range("I26").Value = 13
$ws.Range("J26").Value = 10
$ws.Range("K26").Value = 30
$ws.Range("L26").Value = 44.444444444444
$ws.Range("D15").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Value = 1
$ws.Range("L15").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E27").Value = -100
$ws.Range("A14").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("F27").Formula = "=""0"""
$ws.Range("F27").Copy()
$ws.Range("F27").PasteSpecial(-4163)
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -100
$ws.Range("I27").Value = 16
$ws.Range("J27").Value = 16
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 45.454545454545

$ws.Range("A1").Select() | Out-Null